$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Icam5"
$ws.Cells.Item(2, 3).Value = "Itgal"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 4.589083666666666
$ws.Cells.Item(2, 8).Value = 13.767251
$ws.Cells.Item(2, 9).Value = 0.7170542197587623
$ws.Cells.Item(2, 10).Value = 0.7170542197587624
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 19.43258433333333
$ws.Cells.Item(2, 14).Value = 58.29775300000001
$ws.Cells.Item(2, 15).Value = 0.5025299392457537
$ws.Cells.Item(2, 16).Value = 0.5025299392457538
$ws.Cells.Item(2, 17).Value = 89.17775536522255
$ws.Cells.Item(2, 18).Value = 802.5997982870031
$ws.Cells.Item(2, 19).Value = 0.3603412134912822
$ws.Cells.Item(2, 20).Value = 0.3603412134912823

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Icam5"
$ws.Cells.Item(3, 3).Value = "Itgal"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 4.589083666666666
$ws.Cells.Item(3, 8).Value = 13.767251
$ws.Cells.Item(3, 9).Value = 0.7170542197587623
$ws.Cells.Item(3, 10).Value = 0.7170542197587624
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.3811586666666666
$ws.Cells.Item(3, 14).Value = 1.143476
$ws.Cells.Item(3, 15).Value = 0.009856828012032942
$ws.Cells.Item(3, 16).Value = 0.009856828012032942
$ws.Cells.Item(3, 17).Value = 1.749169011608444
$ws.Cells.Item(3, 18).Value = 15.742521104476
$ws.Cells.Item(3, 19).Value = 0.007067880119464593
$ws.Cells.Item(3, 20).Value = 0.007067880119464594

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Icam5"
$ws.Cells.Item(4, 3).Value = "Itgal"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 4.589083666666666
$ws.Cells.Item(4, 8).Value = 13.767251
$ws.Cells.Item(4, 9).Value = 0.7170542197587623
$ws.Cells.Item(4, 10).Value = 0.7170542197587624
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 18.712703
$ws.Cells.Item(4, 14).Value = 56.138109
$ws.Cells.Item(4, 15).Value = 0.4839136854063913
$ws.Cells.Item(4, 16).Value = 0.4839136854063913
$ws.Cells.Item(4, 17).Value = 85.87415969648433
$ws.Cells.Item(4, 18).Value = 772.867437268359
$ws.Cells.Item(4, 19).Value = 0.3469923501196671
$ws.Cells.Item(4, 20).Value = 0.3469923501196671

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Icam5"
$ws.Cells.Item(5, 3).Value = "Itgal"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4.589083666666666
$ws.Cells.Item(5, 8).Value = 13.767251
$ws.Cells.Item(5, 9).Value = 0.7170542197587623
$ws.Cells.Item(5, 10).Value = 0.7170542197587624
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1430596666666667
$ws.Cells.Item(5, 14).Value = 0.429179
$ws.Cells.Item(5, 15).Value = 0.003699547335821903
$ws.Cells.Item(5, 16).Value = 0.003699547335821903
$ws.Cells.Item(5, 17).Value = 0.6565127796587777
$ws.Cells.Item(5, 18).Value = 5.908615016929
$ws.Cells.Item(5, 19).Value = 0.002652776028348382
$ws.Cells.Item(5, 20).Value = 0.002652776028348383

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Icam5"
$ws.Cells.Item(6, 3).Value = "Itgal"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.495525666666667
$ws.Cells.Item(6, 8).Value = 4.486577
$ws.Cells.Item(6, 9).Value = 0.2336791106752255
$ws.Cells.Item(6, 10).Value = 0.2336791106752255
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 19.43258433333333
$ws.Cells.Item(6, 14).Value = 58.29775300000001
$ws.Cells.Item(6, 15).Value = 0.5025299392457537
$ws.Cells.Item(6, 16).Value = 0.5025299392457538
$ws.Cells.Item(6, 17).Value = 29.06192864016456
$ws.Cells.Item(6, 18).Value = 261.5573577614811
$ws.Cells.Item(6, 19).Value = 0.1174307492906228
$ws.Cells.Item(6, 20).Value = 0.1174307492906229

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Icam5"
$ws.Cells.Item(7, 3).Value = "Itgal"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.495525666666667
$ws.Cells.Item(7, 8).Value = 4.486577
$ws.Cells.Item(7, 9).Value = 0.2336791106752255
$ws.Cells.Item(7, 10).Value = 0.2336791106752255
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.3811586666666666
$ws.Cells.Item(7, 14).Value = 1.143476
$ws.Cells.Item(7, 15).Value = 0.009856828012032942
$ws.Cells.Item(7, 16).Value = 0.009856828012032942
$ws.Cells.Item(7, 17).Value = 0.5700325690724445
$ws.Cells.Item(7, 18).Value = 5.130293121652
$ws.Cells.Item(7, 19).Value = 0.002303334803930509
$ws.Cells.Item(7, 20).Value = 0.002303334803930509

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Icam5"
$ws.Cells.Item(8, 3).Value = "Itgal"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.495525666666667
$ws.Cells.Item(8, 8).Value = 4.486577
$ws.Cells.Item(8, 9).Value = 0.2336791106752255
$ws.Cells.Item(8, 10).Value = 0.2336791106752255
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 18.712703
$ws.Cells.Item(8, 14).Value = 56.138109
$ws.Cells.Item(8, 15).Value = 0.4839136854063913
$ws.Cells.Item(8, 16).Value = 0.4839136854063913
$ws.Cells.Item(8, 17).Value = 27.98532762921034
$ws.Cells.Item(8, 18).Value = 251.867948662893
$ws.Cells.Item(8, 19).Value = 0.1130805196493364
$ws.Cells.Item(8, 20).Value = 0.1130805196493364

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Icam5"
$ws.Cells.Item(9, 3).Value = "Itgal"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.495525666666667
$ws.Cells.Item(9, 8).Value = 4.486577
$ws.Cells.Item(9, 9).Value = 0.2336791106752255
$ws.Cells.Item(9, 10).Value = 0.2336791106752255
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1430596666666667
$ws.Cells.Item(9, 14).Value = 0.429179
$ws.Cells.Item(9, 15).Value = 0.003699547335821903
$ws.Cells.Item(9, 16).Value = 0.003699547335821903
$ws.Cells.Item(9, 17).Value = 0.2139494033647778
$ws.Cells.Item(9, 18).Value = 1.925544630283
$ws.Cells.Item(9, 19).Value = 0.000864506931335762
$ws.Cells.Item(9, 20).Value = 0.0008645069313357622

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Icam5"
$ws.Cells.Item(10, 3).Value = "Itgal"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.3153023333333333
$ws.Cells.Item(10, 8).Value = 0.9459069999999999
$ws.Cells.Item(10, 9).Value = 0.04926666956601224
$ws.Cells.Item(10, 10).Value = 0.04926666956601224
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 19.43258433333333
$ws.Cells.Item(10, 14).Value = 58.29775300000001
$ws.Cells.Item(10, 15).Value = 0.5025299392457537
$ws.Cells.Item(10, 16).Value = 0.5025299392457538
$ws.Cells.Item(10, 17).Value = 6.127139182996777
$ws.Cells.Item(10, 18).Value = 55.144252646971
$ws.Cells.Item(10, 19).Value = 0.02475797646384875
$ws.Cells.Item(10, 20).Value = 0.02475797646384876

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Icam5"
$ws.Cells.Item(11, 3).Value = "Itgal"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.3153023333333333
$ws.Cells.Item(11, 8).Value = 0.9459069999999999
$ws.Cells.Item(11, 9).Value = 0.04926666956601224
$ws.Cells.Item(11, 10).Value = 0.04926666956601224
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3811586666666666
$ws.Cells.Item(11, 14).Value = 1.143476
$ws.Cells.Item(11, 15).Value = 0.009856828012032942
$ws.Cells.Item(11, 16).Value = 0.009856828012032942
$ws.Cells.Item(11, 17).Value = 0.1201802169702222
$ws.Cells.Item(11, 18).Value = 1.081621952732
$ws.Cells.Item(11, 19).Value = 0.0004856130886378403
$ws.Cells.Item(11, 20).Value = 0.0004856130886378403

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Icam5"
$ws.Cells.Item(12, 3).Value = "Itgal"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.3153023333333333
$ws.Cells.Item(12, 8).Value = 0.9459069999999999
$ws.Cells.Item(12, 9).Value = 0.04926666956601224
$ws.Cells.Item(12, 10).Value = 0.04926666956601224
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 18.712703
$ws.Cells.Item(12, 14).Value = 56.138109
$ws.Cells.Item(12, 15).Value = 0.4839136854063913
$ws.Cells.Item(12, 16).Value = 0.4839136854063913
$ws.Cells.Item(12, 17).Value = 5.900158918873666
$ws.Cells.Item(12, 18).Value = 53.10143026986299
$ws.Cells.Item(12, 19).Value = 0.02384081563738788
$ws.Cells.Item(12, 20).Value = 0.02384081563738788

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Icam5"
$ws.Cells.Item(13, 3).Value = "Itgal"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.3153023333333333
$ws.Cells.Item(13, 8).Value = 0.9459069999999999
$ws.Cells.Item(13, 9).Value = 0.04926666956601224
$ws.Cells.Item(13, 10).Value = 0.04926666956601224
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.1430596666666667
$ws.Cells.Item(13, 14).Value = 0.429179
$ws.Cells.Item(13, 15).Value = 0.003699547335821903
$ws.Cells.Item(13, 16).Value = 0.003699547335821903
$ws.Cells.Item(13, 17).Value = 0.04510704670588889
$ws.Cells.Item(13, 18).Value = 0.4059634203529999
$ws.Cells.Item(13, 19).Value = 0.0001822643761377586
$ws.Cells.Item(13, 20).Value = 0.0001822643761377586
